# Resumen_Partidos_AL_2024 - add "Puntos"/standings columns (M:S) to the sheet.
# Mirrors the commit "depuracion de codigo y desarrollo de sidebar":
# adds 7 new computed columns (Posicion resultante, Importancia, Puntos1,
# Puntos2, PuntosDif1, PuntosDif2, Puntos jugados) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("M", "N", "O", "P", "Q", "R", "S")
$headers = @(
    "Posicion resultante",
    "Importancia",
    "Puntos1",
    "Puntos2",
    "PuntosDif1",
    "PuntosDif2",
    "Puntos jugados"
)

# --- header row (row 1) ---
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# --- data rows (2-24): row number followed by M,N,O,P,Q,R,S values ---
$rowsData = @(
    "2|4,3,3,3,0,0,3",
    "3|3,2,6,6,0,0,6",
    "4|7,5,9,9,3,3,9",
    "5|4,4,12,10,3,1,12",
    "6|4,3,13,13,1,1,15",
    "7|4,3,16,16,4,4,18",
    "8|5,5,19,17,7,5,21",
    "9|7,4,20,19,8,7,24",
    "10|6,4,22,21,7,6,27",
    "11|3,5,25,24,7,6,30",
    "12|3,4,28,27,7,6,33",
    "13|3,3,31,30,7,6,36",
    "14|4,5,33,31,9,7,39",
    "15|4,2,34,33,7,6,42",
    "16|4,1,36,34,6,4,45",
    "17|4,2,37,37,4,4,48",
    "18|4,2,40,40,7,7,51",
    "19|2,3,3,1,2,0,3",
    "20|4,4,4,3,3,2,6",
    "21|4,4,5,4,3,2,9",
    "22|4,5,8,5,5,2,12",
    "23|4,5,11,5,4,1,15",
    "24|4,5,14,6,10,2,18"
)

foreach ($rd in $rowsData) {
    $parts = $rd.Split("|")
    $rowNum = [int]$parts[0]
    $vals = $parts[1].Split(",")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rowNum).Value = [int]$vals[$i]
    }
}

# --- match the header formatting (bold + border + centered) used by A1:L1 ---
$ws.Range("L1").Copy()
$ws.Range("M1:S1").PasteSpecial(-4122)

# --- column widths for the new columns (best-fit like A:L) ---
$ws.Range("M1").ColumnWidth = 16.166666666666668
$ws.Range("N1").ColumnWidth = 10.333333333333334
$ws.Range("O1").ColumnWidth = 7
$ws.Range("P1").ColumnWidth = 7
$ws.Range("Q1").ColumnWidth = 9.333333333333334
$ws.Range("R1").ColumnWidth = 9.333333333333334
$ws.Range("S1").ColumnWidth = 13.166666666666666

# --- move the active selection like the saved file shows (T1) ---
$ws.Range("T1").Select() | Out-Null
